$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.519.89"
$ws.Range("E2").Value = "  +1.02%  "
Set-TextValue $ws.Range("D3") "3.351.10"
$ws.Range("E3").Value = "  +1.51%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.28%  "
Set-TextValue $ws.Range("D5") "189.94"
$ws.Range("E5").Value = "  +5.10%  "
Set-TextValue $ws.Range("D6") "559.69"
$ws.Range("E6").Value = "  +0.34%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.18%  "
Set-TextValue $ws.Range("D8") "3.346.95"
$ws.Range("E8").Value = "  +1.64%  "
Set-TextValue $ws.Range("D9") "0.583"
$ws.Range("E9").Value = "  -1.10%  "
Set-TextValue $ws.Range("D10") "0.183"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("E11").Value = "  -0.29%  "
Set-TextValue $ws.Range("D12") "47.05"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("E13").Value = "  +2.44%  "
Set-TextValue $ws.Range("D14") "8.72"
$ws.Range("E14").Value = "  +1.99%  "
Set-TextValue $ws.Range("D15") "3.883.67"
$ws.Range("E15").Value = "  +1.02%  "
Set-TextValue $ws.Range("D16") "602.75"
$ws.Range("E16").Value = "  -4.91%  "
Set-TextValue $ws.Range("D17") "66.531.30"
$ws.Range("E17").Value = "  +0.90%  "
Set-TextValue $ws.Range("D18") "18.06"
$ws.Range("E18").Value = "  +1.10%  "
Set-TextValue $ws.Range("D19") "3.357.86"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("E20").Value = "  +1.15%  "
Set-TextValue $ws.Range("D21") "11.09"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("E22").Value = "  +0.10%  "
Set-TextValue $ws.Range("D23") "18.35"
$ws.Range("E23").Value = "  +4.42%  "
Set-TextValue $ws.Range("D24") "5.05"
$ws.Range("E24").Value = "  +0.52%  "
Set-TextValue $ws.Range("D25") "100.49"
$ws.Range("E25").Value = "  -5.88%  "
Set-TextValue $ws.Range("D26") "4.00"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").Value = "  +2.64%  "
Set-TextValue $ws.Range("D29") "9.61"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  -0.08%  "
Set-TextValue $ws.Range("D31") "30.81"
$ws.Range("E31").Value = "  +0.75%  "
Set-TextValue $ws.Range("D32") "6.73"
$ws.Range("E32").Value = "  +5.88%  "
$ws.Range("E33").Value = "  -1.36%  "
Set-TextValue $ws.Range("D34") "588.00"
$ws.Range("E34").Value = "  +6.68%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  +0.28%  "
Set-TextValue $ws.Range("D37") "3.742.69"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E38").Value = "  +0.04%  "
Set-TextValue $ws.Range("D39") "56.41"
$ws.Range("E39").Value = "  -1.23%  "
Set-TextValue $ws.Range("D40") "3.54"
$ws.Range("E40").Value = "  +4.57%  "
Set-TextValue $ws.Range("D41") "33.96"
$ws.Range("E41").Value = "  +5.32%  "
Set-TextValue $ws.Range("D42") "0.0₃0713"
$ws.Range("E42").Value = "  -0.15%  "
Set-TextValue $ws.Range("D43") "3.25"
$ws.Range("E43").Value = "  -6.96%  "
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +0.48%  "
Set-TextValue $ws.Range("D47") "3.38"
$ws.Range("E47").Value = "  +6.06%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E49").Value = "  +0.20%  "
Set-TextValue $ws.Range("D50") "2.60"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  -0.25%  "
